# Product backlog.xlsx - "Atualização das atividade do Backlog final"
#
# The sprint is over: every remaining item on the Product Backlog sheet
# (previously a mix of "Pendente" and "Backlog" statuses) is now marked
# "Concluído" in the Status column (L2:L55). Because no cell references
# the "Pendente" / "Backlog" shared strings any more, they fall out of
# the shared-string table on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product Backlog")

# Column L ("Status") runs from row 2 through row 55 (dimension B1:L59).
for ($r = 2; $r -le 55; $r++) {
    $ws.Range("L$r").Value = "Concluído"
}

# Reflect the author's final view position/selection on the sheet: the
# frozen header pane stays frozen at row 1, the scrollable body pane is
# positioned so row 53 is at the top, and the active selection moves to I54.
$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A2").Select() | Out-Null
$win.FreezePanes = $true

$panes = $win.Panes
$bodyPane = $panes.Item($panes.Count)
$bodyPane.ScrollRow = 53
$bodyPane.ScrollColumn = 1

$ws.Range("I54").Select() | Out-Null
